$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.846.16"
$ws.Range("E2").Value = "  +4.33%  "
$ws.Range("D3").Value = "2.275.71"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.06"
$ws.Range("E5").Value = "  +4.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.74"
$ws.Range("E6").Value = "  +5.31%  "
$ws.Range("E7").Value = "  +3.74%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.63"
$ws.Range("E10").Value = "  +5.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.76"
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("D15").Value = "2.627.43"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.28"
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").Value = "2.274.70"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("E18").Value = "  +3.64%  "
$ws.Range("D19").Value = "41.792.76"
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("E20").Value = "  +8.75%  "
$ws.Range("D21").Value = "0.0₃0909"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.27"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.63"
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("E25").Value = "  +5.28%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  +4.48%  "
$ws.Range("E28").Value = "  +4.77%  "
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.07"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.53"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.07"
$ws.Range("E32").Value = "  +6.84%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0749"
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.85"
$ws.Range("E38").Value = "  +7.07%  "
$ws.Range("E39").Value = "  +5.39%  "
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("E41").Value = "  +3.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.94"
$ws.Range("E42").Value = "  +5.15%  "
$ws.Range("D43").Value = "2.077.60"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.64"
$ws.Range("E44").Value = "  +3.90%  "
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.31"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  +7.61%  "
$ws.Range("E48").Value = "  +4.76%  "
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.97"
$ws.Range("E50").Value = "  +7.57%  "
$ws.Range("E51").Value = "  +3.65%  "
